# edit.ps1
# Applies the diff:
#  1. Inserts two new numbered (Issue/Resolution) paragraphs about
#     "Categorizing the data based on movie categories" plus a trailing
#     empty (ind left=1440) paragraph, right before the blank paragraph
#     that precedes the "9. Conclusion" heading.
#  2. Splits the run in the Conclusion paragraph so that a
#     <w:lastRenderedPageBreak/> is inserted before "of trailers...".

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: insert the new "Issue"/"Resolution" list paragraphs plus the
# blank indented paragraph, right before the empty paragraph that comes
# immediately after the "Resolution: Used Flutter's ..." bullet.
# ---------------------------------------------------------------------

$anchorRange = $d.Content
$anchorRange.Find.Execute("Resolution: Used Flutter") | Out-Null
$anchorPara = $anchorRange.Paragraphs(1)
$blankPara = $anchorPara.Next()

# Create three fresh empty paragraphs right before $blankPara (which is
# the existing empty "<w:p><w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>"
# paragraph that precedes the Conclusion heading).
$blankPara.Range.InsertParagraphBefore()
$blankPara.Range.InsertParagraphBefore()
$blankPara.Range.InsertParagraphBefore()

$newPara1 = $anchorPara.Next()
$newPara2 = $newPara1.Next()
$newPara3 = $newPara2.Next()

$issueXml = "<w:p $wNs>" +
              "<w:pPr>" +
                "<w:numPr>" +
                  "<w:ilvl w:val=`"0`"/>" +
                  "<w:numId w:val=`"7`"/>" +
                "</w:numPr>" +
              "</w:pPr>" +
              "<w:r>" +
                "<w:rPr>" +
                  "<w:b/>" +
                  "<w:bCs/>" +
                "</w:rPr>" +
                "<w:t>Issue</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:t xml:space=`"preserve`">: </w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:t>Categorizing the data based on movie categories</w:t>" +
              "</w:r>" +
              "<w:r>" +
                "<w:t>.</w:t>" +
              "</w:r>" +
            "</w:p>"
$newPara1.Range.InsertXML($issueXml) | Out-Null

$resolutionXml = "<w:p $wNs>" +
                    "<w:pPr>" +
                      "<w:numPr>" +
                        "<w:ilvl w:val=`"1`"/>" +
                        "<w:numId w:val=`"7`"/>" +
                      "</w:numPr>" +
                    "</w:pPr>" +
                    "<w:r>" +
                      "<w:t>Resolution:</w:t>" +
                    "</w:r>" +
                    "<w:r>" +
                      "<w:t xml:space=`"preserve`"> Implement the use of dynamic maps to store each value of data into </w:t>" +
                    "</w:r>" +
                    "<w:proofErr w:type=`"gramStart`"/>" +
                    "<w:r>" +
                      "<w:t>specified</w:t>" +
                    "</w:r>" +
                    "<w:proofErr w:type=`"gramEnd`"/>" +
                    "<w:r>" +
                      "<w:t xml:space=`"preserve`"> key which is the category of the movie</w:t>" +
                    "</w:r>" +
                    "<w:r>" +
                      "<w:t>.</w:t>" +
                    "</w:r>" +
                  "</w:p>"
$newPara2.Range.InsertXML($resolutionXml) | Out-Null

$blankIndentXml = "<w:p $wNs>" +
                     "<w:pPr>" +
                       "<w:ind w:left=`"1440`"/>" +
                     "</w:pPr>" +
                   "</w:p>"
$newPara3.Range.InsertXML($blankIndentXml) | Out-Null

# ---------------------------------------------------------------------
# Part 2: split the run containing "The integration of trailers ..." so
# that a <w:lastRenderedPageBreak/> appears right before "of trailers".
# The whole containing paragraph is rebuilt (with the run split in two)
# and used to replace the paragraph's full range in a single shot -
# this keeps everything inside the same <w:p> instead of the runtime's
# InsertXML spinning off a brand-new paragraph.
# ---------------------------------------------------------------------

$concRange = $d.Content
$concRange.Find.Execute("The integration of trailers and real-time data significantly enhances user engagement, while the MVVM architecture and Provider state management ensure that the app is scalable and maintainable.") | Out-Null
$concPara = $concRange.Paragraphs(1)
$fullParaRange = $concPara.Range

$newConclusionParaXml = "<w:p $wNs>" +
    "<w:r><w:t xml:space=`"preserve`">The </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/>" +
    "<w:r><w:t>Moviepedia</w:t></w:r>" +
    "<w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> app successfully achieved its primary goals, providing users with a comprehensive, user-friendly, and visually appealing platform for exploring movies. The integration </w:t></w:r>" +
    "<w:r>" +
      "<w:lastRenderedPageBreak/>" +
      "<w:t>of trailers and real-time data significantly enhances user engagement, while the MVVM architecture and Provider state management ensure that the app is scalable and maintainable.</w:t>" +
    "</w:r>" +
  "</w:p>"

$fullParaRange.InsertXML($newConclusionParaXml) | Out-Null

Write-Host "Edit applied successfully."
